$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(33191000,29778000,42244000,45441000,48709000,52066000,55485000,62200000,69045000,85034000,102048000,153106000,211304000,291561000,402326000,555186000,35532000,30892000,50308000,54080000,57853000,61847000,64954000,73489000,82383000,101404000,121685000,182507000,251898000,347640000,479750000,662028000,40684000,35370000,63673000,68214000,72818000,77372000,83436000,88107000,101342000,114037000,157692000,212894000,298967000,415869000,582216000,815102000)

for ($row = 50; $row -le 289; $row++) {
    $idx = ($row - 50) % 48
    $ws.Cells.Item($row, 7).Value = $newValues[$idx]
}

for ($row = 50; $row -le 145; $row++) {
    $ws.Rows.Item($row).RowHeight = 13.8
}

for ($row = 146; $row -le 289; $row++) {
    $ws.Rows.Item($row).RowHeight = 66.55
}

$ws.Application.ActiveWindow.ScrollRow = 298
$ws.Range("G301").Select()
